# Kas HIMA - update 15 feb 2023
# Adds the two new ledger entries for the valentine-related transactions
# (expense + offering) that occurred on 13 Feb 2023 and 15 Feb 2023, and
# moves the active selection/scroll position to where the user left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$currencyFormat = '_-"Rp"* #,##0.00_-;\-"Rp"* #,##0.00_-;_-"Rp"* "-"??_-;_-@_-'

# --- Row 14: pengeluaran belanja buat valentine (13 Feb 2023) -----------
$ws.Range("A14").Value = 44970
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 436000
$ws.Range("D14").Formula = "=D13+B14-C14"
$ws.Range("B14:D14").NumberFormat = $currencyFormat
$ws.Range("E14").Value = "pengeluaran belanja buat valentine"

# --- Row 15: uang persembahan - reguler - valentine (15 Feb 2023) -------
$ws.Range("A15").Value = 44972
$ws.Range("B15").Value = 95500
$ws.Range("C15").Value = 0
$ws.Range("D15").Formula = "=D14+B15-C15"
$ws.Range("B15:D15").NumberFormat = $currencyFormat
$ws.Range("E15").Value = "uang persembahan - reguler - valentine"

# --- View state: scroll right a column and park the selection on C21 ----
$win = $wb.Windows.Item(1)
$win.ScrollColumn = 2
$ws.Range("C21").Select()
